$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting so numeric-looking strings (e.g. "2.1", "1") are
# stored as text, matching the workbook's existing inline-string cells
# instead of being auto-coerced to numbers.
$textCells = @("K2","L2","K3","L3","AO3","K4","L4","AO4","K5","L5","K6","L6","AO6","K7","L7","AO7")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2: strip the bracket wrapper from the error message and add the new
# Statement ID column value.
$ws.Range("K2").Value = "The 'Encoded Statement' field does not contain IG Script-encoded content."
$ws.Range("L2").Value = "1"

# Row 3: "[OK]" -> "OK"; Statement ID offset by the new row; linkage updated.
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "2.1"
$ws.Range("AO3").Value = "[OR].Bdir.[2.2]"

# Row 4
$ws.Range("K4").Value = "OK"
$ws.Range("L4").Value = "2.2"
$ws.Range("AO4").Value = "[OR].Bdir.[2.1]"

# Row 5
$ws.Range("K5").Value = "OK"
$ws.Range("L5").Value = "3"

# Row 6
$ws.Range("K6").Value = "OK"
$ws.Range("L6").Value = "4.1"
$ws.Range("AO6").Value = "[OR].Bdir.[4.2]"

# Row 7
$ws.Range("K7").Value = "OK"
$ws.Range("L7").Value = "4.2"
$ws.Range("AO7").Value = "[OR].Bdir.[4.1]"
